# edit.ps1 - applies the commit's changes to the presentation:
#   1) Swap the deck's theme colour scheme from "Integral" (Red Violet)
#      to the stock "Office" colour scheme (the two theme parts'
#      clrScheme content is swapped in the canonical OOXML; the
#      colour-scheme half of that is reachable through the PowerPoint
#      object model via Slide.ThemeColorScheme, which edits the theme
#      part shared by the whole deck).
#   2) Re-point the three data tables (slides 14, 15, 16) from the
#      custom table style {C3A4D8BA-4FC4-49BE-AD23-F80122E59843} to the
#      built-in style {E8DFB77F-A7F8-47DA-84F6-17C2AB90CAF9}.

$p = $ppt.ActivePresentation

# --- 1) Theme colours: Integral/"Red Violet" -> stock "Office" scheme ---
# Order is dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink (indices 1-12),
# values are packed 0xBBGGRR the way VBA's RGB()/ColorFormat.RGB works.
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}

# --- 2) Table styles on the three table slides ---
$newStyleId = "{E8DFB77F-A7F8-47DA-84F6-17C2AB90CAF9}"
foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId, $true)
        }
    }
}
